$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.926.58"
$ws.Range("E2").Value = "  +6.11%  "
$ws.Range("D3").Value = "1.879.14"
$ws.Range("E3").Value = "  +5.38%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.55"
$ws.Range("E5").Value = "  +1.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4977"
$ws.Range("E7").Value = "  +1.24%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "45.91"
$ws.Range("E8").Value = "  +8.81%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2851"
$ws.Range("E9").Value = "  +6.83%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06526"
$ws.Range("E10").Value = "  +4.44%  "
$ws.Range("D11").Value = "1.879.51"
$ws.Range("E11").Value = "  +5.83%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "17.05"
$ws.Range("E12").Value = "  +3.46%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07234"
$ws.Range("E13").Value = "  +3.14%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6638"
$ws.Range("E14").Value = "  +6.00%  "
$ws.Range("E15").Value = "  +6.49%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.798"
$ws.Range("E16").Value = "  +3.53%  "
$ws.Range("D17").Value = "29.951.88"
$ws.Range("E17").Value = "  +6.33%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.000"
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.86"
$ws.Range("E19").Value = "  +6.72%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007496"
$ws.Range("E20").Value = "  +3.72%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("D22").Value = "2.123.09"
$ws.Range("E22").Value = "  +5.76%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.739"
$ws.Range("E23").Value = "  +3.91%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.536"
$ws.Range("E24").Value = "  +5.80%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.000"
$ws.Range("E25").Value = "  +3.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "145.56"
$ws.Range("E26").Value = "  +3.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "134.94"
$ws.Range("E27").Value = "  +23.72%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.66"
$ws.Range("E28").Value = "  +5.65%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.950"
$ws.Range("E29").Value = "  +4.92%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.383"
$ws.Range("E30").Value = "  -0.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.179"
$ws.Range("E31").Value = "  +0.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08600"
$ws.Range("E32").Value = "  +4.32%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.860"
$ws.Range("E33").Value = "  +2.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05111"
$ws.Range("E34").Value = "  +4.64%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.127"
$ws.Range("E35").Value = "  +5.08%  "
$ws.Range("B36").Value = "Frax"
$ws.Range("C36").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.000"
$ws.Range("E36").Value = "  +0.09%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6845"
$ws.Range("E37").Value = "  +5.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.710"
$ws.Range("E38").Value = "  +3.73%  "
$ws.Range("E39").Value = "  +12.48%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.743"
$ws.Range("E40").Value = "  +5.98%  "
$ws.Range("E41").Value = "  +1.25%  "
$ws.Range("E42").Value = "  +5.19%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.092"
$ws.Range("E43").Value = "  +2.83%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "104.29"
$ws.Range("E44").Value = "  +4.85%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.001"
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4207"
$ws.Range("E46").Value = "  +5.81%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.425"
$ws.Range("E47").Value = "  +3.35%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1249"
$ws.Range("E48").Value = "  +3.77%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05637"
$ws.Range("E49").Value = "  +3.87%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "32.30"
$ws.Range("E50").Value = "  +5.52%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.232"
$ws.Range("E51").Value = "  +2.97%  "
